# 自动更新Excel文件 - 2025-10-22 23:11:00
# Recompute "剩余" (E, remaining days) for each row against the new
# reference date (2025-10-23). Remaining = 总天(D) - (today - 开始时间(F)).
# When that would drop to zero or below, the cycle restarts: remaining
# resets to the full duration (D) and the start date (F) resets to today.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference "today" used to recompute the remaining-day counters.
$today = Get-Date -Year 2025 -Month 10 -Day 23 -Hour 0 -Minute 0 -Second 0
$todayOA = $today.ToOADate()
$todayNum = [int]$today.ToString("yyyyMMdd")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells($r, 4).Value2
    $eVal = $ws.Cells($r, 5).Value2
    $fVal = $ws.Cells($r, 6).Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string][int]$fVal
    if ($fStr.Length -ne 8) {
        # Malformed start date (data-entry typo) - leave the row untouched.
        continue
    }

    $y = [int]$fStr.Substring(0, 4)
    $mo = [int]$fStr.Substring(4, 2)
    $da = [int]$fStr.Substring(6, 2)
    $fDate = Get-Date -Year $y -Month $mo -Day $da -Hour 0 -Minute 0 -Second 0
    $fOA = $fDate.ToOADate()

    $elapsed = $todayOA - $fOA
    $newE = $dVal - $elapsed

    if ($newE -le 0) {
        # Cycle completed - restart it as of today.
        $newE = $dVal
        $newF = $todayNum
    } else {
        $newF = $fVal
    }

    if ($newE -ne $eVal) {
        $ws.Cells($r, 5).Value2 = $newE
    }
    if ($newF -ne $fVal) {
        $ws.Cells($r, 6).Value2 = $newF
    }
}
